# This script updates the simulated transition-probability matrix on Sheet1.
# More games were simulated (see commit message), which changed the tallied
# outcome counts underlying each row's probabilities; the recomputed
# probabilities (updated counts / updated row totals) are written below.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 0.2520661157024793
$ws.Range("C2").Value = 0.4586776859504132
$ws.Range("J2").Value = 0.03305785123966942
$ws.Range("P2").Value = 0.1818181818181818
$ws.Range("S2").Value = 0.0743801652892562

# Row 3
$ws.Range("B3").Value = 0.008849557522123894
$ws.Range("C3").Value = 0.02654867256637168
$ws.Range("J3").Value = 0.03539823008849557
$ws.Range("P3").Value = 0.7876106194690266
$ws.Range("S3").Value = 0.1415929203539823

# Row 4
$ws.Range("J4").Value = 0.08108108108108109
$ws.Range("P4").Value = 0.7297297297297297
$ws.Range("S4").Value = 0.1891891891891892

# Row 6
$ws.Range("B6").Value = 0.05583756345177665
$ws.Range("D6").Value = 0.01015228426395939
$ws.Range("F6").Value = 0.06598984771573604
$ws.Range("J6").Value = 0.233502538071066
$ws.Range("O6").Value = 0.01522842639593909
$ws.Range("Q6").Value = 0.1725888324873096
$ws.Range("R6").Value = 0.08121827411167512
$ws.Range("S6").Value = 0.3654822335025381

# Row 7
$ws.Range("B7").Value = 0.07017543859649122
$ws.Range("D7").Value = 0.05847953216374269
$ws.Range("F7").Value = 0.04678362573099415
$ws.Range("J7").Value = 0.1286549707602339
$ws.Range("O7").Value = 0.005847953216374269
$ws.Range("Q7").Value = 0.1637426900584795
$ws.Range("R7").Value = 0.07602339181286549
$ws.Range("S7").Value = 0.4502923976608187

# Row 8
$ws.Range("B8").Value = 0.06772009029345373
$ws.Range("D8").Value = 0.009029345372460496
$ws.Range("F8").Value = 0.05869074492099323
$ws.Range("J8").Value = 0.1422121896162528
$ws.Range("O8").Value = 0.02257336343115124
$ws.Range("Q8").Value = 0.1625282167042889
$ws.Range("R8").Value = 0.1060948081264108
$ws.Range("S8").Value = 0.4311512415349887

# Row 9
$ws.Range("B9").Value = 0.0778688524590164
$ws.Range("D9").Value = 0.01229508196721311
$ws.Range("F9").Value = 0.04918032786885246
$ws.Range("J9").Value = 0.110655737704918
$ws.Range("O9").Value = 0.01229508196721311
$ws.Range("Q9").Value = 0.2213114754098361
$ws.Range("R9").Value = 0.0860655737704918
$ws.Range("S9").Value = 0.430327868852459

# Row 10
$ws.Range("B10").Value = 0.09052808046940486
$ws.Range("D10").Value = 0.01676445934618609
$ws.Range("E10").Value = 0.0008382229673093043
$ws.Range("F10").Value = 0.06454316848281642
$ws.Range("J10").Value = 0.124056999161777
$ws.Range("O10").Value = 0.01508801341156748
$ws.Range("Q10").Value = 0.2321877619446773
$ws.Range("R10").Value = 0.09388097233864208
$ws.Range("S10").Value = 0.3621123218776194

# Row 11
$ws.Range("G11").Value = 0.1145374449339207
$ws.Range("J11").Value = 0.08370044052863436
$ws.Range("K11").Value = 0.1674008810572687
$ws.Range("L11").Value = 0.6167400881057269
$ws.Range("S11").Value = 0.01762114537444934

# Row 12
$ws.Range("G12").Value = 0.8014184397163121
$ws.Range("J12").Value = 0.1276595744680851
$ws.Range("K12").Value = 0.007092198581560284
$ws.Range("L12").Value = 0.02836879432624113
$ws.Range("S12").Value = 0.03546099290780142

# Row 13
$ws.Range("G13").Value = 0.7446808510638298
$ws.Range("J13").Value = 0.2127659574468085
$ws.Range("S13").Value = 0.0425531914893617

# Row 15
$ws.Range("F15").Value = 0.03571428571428571
$ws.Range("H15").Value = 0.1830357142857143
$ws.Range("I15").Value = 0.06696428571428571
$ws.Range("J15").Value = 0.3705357142857143
$ws.Range("K15").Value = 0.04910714285714286
$ws.Range("M15").Value = 0.01339285714285714
$ws.Range("O15").Value = 0.04464285714285714
$ws.Range("S15").Value = 0.2366071428571428

# Row 16
$ws.Range("F16").Value = 0.006451612903225806
$ws.Range("H16").Value = 0.2129032258064516
$ws.Range("I16").Value = 0.06451612903225806
$ws.Range("J16").Value = 0.4064516129032258
$ws.Range("K16").Value = 0.1096774193548387
$ws.Range("M16").Value = 0.01935483870967742
$ws.Range("O16").Value = 0.06451612903225806
$ws.Range("S16").Value = 0.1161290322580645

# Row 17
$ws.Range("F17").Value = 0.01952277657266811
$ws.Range("H17").Value = 0.1496746203904555
$ws.Range("I17").Value = 0.1127982646420824
$ws.Range("J17").Value = 0.4229934924078091
$ws.Range("K17").Value = 0.09978308026030369
$ws.Range("M17").Value = 0.01301518438177874
$ws.Range("O17").Value = 0.07158351409978309
$ws.Range("S17").Value = 0.1106290672451193

# Row 18
$ws.Range("F18").Value = 0.02870813397129187
$ws.Range("H18").Value = 0.1626794258373206
$ws.Range("I18").Value = 0.1100478468899522
$ws.Range("J18").Value = 0.4019138755980861
$ws.Range("K18").Value = 0.09090909090909091
$ws.Range("M18").Value = 0.01913875598086124
$ws.Range("O18").Value = 0.1004784688995215
$ws.Range("S18").Value = 0.0861244019138756

# Row 19
$ws.Range("F19").Value = 0.01167639699749791
$ws.Range("H19").Value = 0.2243536280233528
$ws.Range("I19").Value = 0.1209341117597998
$ws.Range("J19").Value = 0.3427856547122602
$ws.Range("K19").Value = 0.07923269391159299
$ws.Range("M19").Value = 0.02668890742285238
$ws.Range("O19").Value = 0.07172643869891576
$ws.Range("S19").Value = 0.1226021684737281

Write-Host "Updated 107 cells on Sheet1."
